# B1--and-B2-PowerPoint.pptx - Mon, Aug 03, 2020 6:05:05 AM
#
# 1) The table on slide 5 gets switched to a different built-in table style.
# 2) The deck's theme colour scheme (custom "Red Violet"/Integral palette)
#    is swapped for the standard Office theme palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the comparison table on slide 5 -----------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{D00523CD-ACAE-4E1A-9D40-B693A2CA24E9}")

# --- 2. Swap the theme colours for the standard "Office" palette -------
$firstSlide = $p.Slides.Item(1)
$themeColors = $firstSlide.ThemeColorScheme

$themeColors.Item(1).RGB  = 0        # Text/Background - Dark 1  -> 000000
$themeColors.Item(2).RGB  = 16777215 # Text/Background - Light 1 -> FFFFFF
$themeColors.Item(3).RGB  = 6968388  # Text/Background - Dark 2  -> 44546A
$themeColors.Item(4).RGB  = 15132391 # Text/Background - Light 2 -> E7E6E6
$themeColors.Item(5).RGB  = 13998939 # Accent 1 -> 5B9BD5
$themeColors.Item(6).RGB  = 3243501  # Accent 2 -> ED7D31
$themeColors.Item(7).RGB  = 10855845 # Accent 3 -> A5A5A5
$themeColors.Item(8).RGB  = 49407    # Accent 4 -> FFC000
$themeColors.Item(9).RGB  = 12874308 # Accent 5 -> 4472C4
$themeColors.Item(10).RGB = 4697456  # Accent 6 -> 70AD47
$themeColors.Item(11).RGB = 12673797 # Hyperlink -> 0563C1
$themeColors.Item(12).RGB = 7491477  # Followed Hyperlink -> 954F72
